$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that belong to the accidentally-included extra condition
# blocks (category 3 tail, category 5 tail, category 4 tail). Deleting from
# the bottom up keeps the remaining row numbers stable while we work.
$ws.Rows.Item(101).Delete()
$ws.Rows.Item(100).Delete()
$ws.Rows.Item(99).Delete()
$ws.Rows.Item(98).Delete()
$ws.Rows.Item(97).Delete()
$ws.Rows.Item(96).Delete()
$ws.Rows.Item(81).Delete()
$ws.Rows.Item(80).Delete()
$ws.Rows.Item(79).Delete()
$ws.Rows.Item(78).Delete()
$ws.Rows.Item(77).Delete()
$ws.Rows.Item(76).Delete()
$ws.Rows.Item(75).Delete()
$ws.Rows.Item(61).Delete()
$ws.Rows.Item(60).Delete()
$ws.Rows.Item(59).Delete()
$ws.Rows.Item(58).Delete()
$ws.Rows.Item(57).Delete()
$ws.Rows.Item(56).Delete()
$ws.Rows.Item(55).Delete()

$ws.Range("D68").Select()
